$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.231316089630127
$ws.Range("B1").Value = 2.505062341690063
$ws.Range("C1").Value = 4.384113311767578
$ws.Range("D1").Value = 2.557069301605225
$ws.Range("E1").Value = 1.077600836753845
